$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Dll1","Notch1","ECs",3.0,1.0,16.39002933333333,49.170088,0.5551882184054378,0.5551882184054378,3.0,1.0,70.462916,211.388748,0.5276750397950939,0.5276750397950939,1154.889260152203,10394.00334136983,0.2929589652408566,0.2929589652408566),
    @("ECs","Dll1","Notch1","FAPs",3.0,1.0,16.39002933333333,49.170088,0.5551882184054378,0.5551882184054378,3.0,1.0,9.829723666666666,29.489171,0.07361176802536967,0.07361176802536967,161.1094592352275,1449.985133117048,0.04086838634367936,0.04086838634367936),
    @("ECs","Dll1","Notch1","MuSCs",3.0,1.0,16.39002933333333,49.170088,0.5551882184054378,0.5551882184054378,3.0,1.0,42.505498,127.516494,0.3183105613832428,0.3183105613832428,696.6663590479413,6269.997231431472,0.1767222734739973,0.1767222734739973),
    @("ECs","Dll1","Notch1","Resolving-Mac",3.0,1.0,16.39002933333333,49.170088,0.5551882184054378,0.5551882184054378,3.0,1.0,10.73653933333333,32.209618,0.08040263079629371,0.08040263079629371,175.9721946118204,1583.749751506384,0.04463859334690449,0.04463859334690449),
    @("FAPs","Dll1","Notch1","ECs",2.0,0.6666666666666666,0.2120556666666667,0.636167,0.007183074867352934,0.007183074867352935,3.0,1.0,70.462916,211.388748,0.5276750397950939,0.5276750397950939,14.94206062765734,134.478545648916,0.003790329316481598,0.003790329316481599),
    @("FAPs","Dll1","Notch1","FAPs",2.0,0.6666666666666666,0.2120556666666667,0.636167,0.007183074867352934,0.007183074867352935,3.0,1.0,9.829723666666666,29.489171,0.07361176802536967,0.07361176802536967,2.084448605284111,18.760037447557,0.0005287588408444472,0.0005287588408444473),
    @("FAPs","Dll1","Notch1","MuSCs",2.0,0.6666666666666666,0.2120556666666667,0.636167,0.007183074867352934,0.007183074867352935,3.0,1.0,42.505498,127.516494,0.3183105613832428,0.3183105613832428,9.013531715388666,81.121785438498,0.002286448593484975,0.002286448593484975),
    @("FAPs","Dll1","Notch1","Resolving-Mac",2.0,0.6666666666666666,0.2120556666666667,0.636167,0.007183074867352934,0.007183074867352935,3.0,1.0,10.73653933333333,32.209618,0.08040263079629371,0.08040263079629371,2.276744006022889,20.490696054206,0.0005775381165419144,0.0005775381165419145),
    @("MuSCs","Dll1","Notch1","ECs",3.0,1.0,12.89091033333333,38.672731,0.4366606914505164,0.4366606914505165,3.0,1.0,70.462916,211.388748,0.5276750397950939,0.5276750397950939,908.3311319811987,8174.980187830789,0.2304149477381044,0.2304149477381045),
    @("MuSCs","Dll1","Notch1","FAPs",3.0,1.0,12.89091033333333,38.672731,0.4366606914505164,0.4366606914505165,3.0,1.0,9.829723666666666,29.489171,0.07361176802536967,0.07361176802536967,126.7140863884445,1140.426777496001,0.03214336552485293,0.03214336552485295),
    @("MuSCs","Dll1","Notch1","MuSCs",3.0,1.0,12.89091033333333,38.672731,0.4366606914505164,0.4366606914505165,3.0,1.0,42.505498,127.516494,0.3183105613832428,0.3183105613832428,547.9345633916793,4931.411070525113,0.1389937098296088,0.1389937098296089),
    @("MuSCs","Dll1","Notch1","Resolving-Mac",3.0,1.0,12.89091033333333,38.672731,0.4366606914505164,0.4366606914505165,3.0,1.0,10.73653933333333,32.209618,0.08040263079629371,0.08040263079629371,138.4037658363064,1245.633892526758,0.0351086683579502,0.03510866835795021),
    @("Resolving-Mac","Dll1","Notch1","ECs",1.0,0.3333333333333333,0.02857733333333333,0.085732,0.0009680152766929151,0.0009680152766929153,3.0,1.0,70.462916,211.388748,0.5276750397950939,0.5276750397950939,2.013642238170667,18.122780143536,0.0005107974996511928,0.0005107974996511929),
    @("Resolving-Mac","Dll1","Notch1","FAPs",1.0,0.3333333333333333,0.02857733333333333,0.085732,0.0009680152766929151,0.0009680152766929153,3.0,1.0,9.829723666666666,29.489171,0.07361176802536967,0.07361176802536967,0.2809072897968889,2.528165608172,0.00007125731599293291,0.00007125731599293292),
    @("Resolving-Mac","Dll1","Notch1","MuSCs",1.0,0.3333333333333333,0.02857733333333333,0.085732,0.0009680152766929151,0.0009680152766929153,3.0,1.0,42.505498,127.516494,0.3183105613832428,0.3183105613832428,1.214693784845333,10.932244063608,0.0003081294861516769,0.000308129486151677),
    @("Resolving-Mac","Dll1","Notch1","Resolving-Mac",1.0,0.3333333333333333,0.02857733333333333,0.085732,0.0009680152766929151,0.0009680152766929153,3.0,1.0,10.73653933333333,32.209618,0.08040263079629371,0.08040263079629371,0.3068216633751111,2.761394970376,0.00007783097489711256,0.00007783097489711257)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
